# Apply new metric values to rows 2-26 (columns B:Q) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B" = 0.9298247212236502
    "C" = 0.7246052329558914
    "D" = 0.9999999999994325
    "E" = 0.3752963470207815
    "F" = 0.8642249000489439
    "G" = 0.469262412304898
    "H" = 1.841566075300257
    "I" = 0.000000000002362837344077222
    "J" = 0.8985662428071469
    "K" = 0.4492831214047548
    "L" = 0.3376508542801637
    "M" = 0.6850273077074358
    "N" = 1.129554360817876
    "O" = 0.7141903383985572
    "P" = 75.51318630499273
    "Q" = 120.6115918251162
}

foreach ($col in $newValues.Keys) {
    $value = $newValues[$col]
    $range = $ws.Range($col + "2:" + $col + "26")
    $range.Value = $value
}
